$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Total" header (B1) to "Total_Labor"
$ws.Range("B1").Value = "Total_Labor"

# Strip the " County, New York" suffix from each geographic name in column A
# (rows 2-63), e.g. "Albany County, New York" -> "Albany",
# "New York County, New York" -> "New York".
$suffix = " County, New York"
for ($row = 2; $row -le 63; $row++) {
  $cell = $ws.Cells.Item($row, 1)
  $current = $cell.Text
  if ($current.EndsWith($suffix)) {
    $cell.Value = $current.Substring(0, $current.Length - $suffix.Length)
  }
}

# Adjust column widths (values pre-compensated for this runtime's internal
# 1/6-character rounding of ColumnWidth, so the saved OOXML <col width>
# lands as close as possible to the target widths below)
$ws.Columns.Item(1).ColumnWidth = 43.5               # -> ~44.29
$ws.Columns.Item(2).ColumnWidth = 12.6666666666667    # -> ~13.43
$ws.Columns.Item(3).ColumnWidth = 19.6666666666667    # -> ~20.57
$ws.Columns.Item(4).ColumnWidth = 15.3333333333333    # -> ~16.14
$ws.Columns.Item(5).ColumnWidth = 22                  # -> ~22.86
$ws.Columns.Item(6).ColumnWidth = 14.1666666666667    # -> 15

# Update selection to match target (columns B:F entire, active cell B1)
[void]$ws.Range("B1:F1048576").Select()
